$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full data for rows 2-10 (Sending cluster, Ligand, Receptor, Target cluster, then metric columns E..T)
$data = @(
    @("ECs","Ntf3","Ntrk1","ECs",3,1,1.679513333333333,5.038539999999999,0.376631045782902,0.3928053077137587,1,0.3333333333333333,0.0009183333333333334,0.002755,0.005182370534790543,0.005409247358696911,0.001542353077777778,0.0138811777,0.001951841634152659,0.002124781073232777),
    @("ECs","Ntf3","Ntrk1","FAPs",3,1,1.679513333333333,5.038539999999999,0.376631045782902,0.3928053077137587,2,0.6666666666666666,0.153988,0.461964,0.8689904253117887,0.9070335923096405,0.2586248991733333,2.32762409256,0.3272887726605078,0.3562876093339043),
    @("ECs","Ntf3","Ntrk1","MuSCs",3,1,1.679513333333333,5.038539999999999,0.376631045782902,0.3928053077137587,1,0.5,0.022297,0.044594,0.1258272041534207,0.08755716033166246,0.03744810879333333,0.22468865276,0.04739043148824156,0.03439291730662158),
    @("FAPs","Ntf3","Ntrk1","ECs",3,1,2.228940333333334,6.686821,0.4998401096732527,0.5213055330575571,1,0.3333333333333333,0.0009183333333333334,0.002755,0.005182370534790543,0.005409247358696911,0.002046910206111112,0.018422191855,0.002590356656477138,0.002819870577765676),
    @("FAPs","Ntf3","Ntrk1","FAPs",3,1,2.228940333333334,6.686821,0.4998401096732527,0.5213055330575571,2,0.6666666666666666,0.153988,0.461964,0.8689904253117887,0.9070335923096405,0.3432300640493333,3.089070576444,0.434356269492851,0.4728416303400881),
    @("FAPs","Ntf3","Ntrk1","MuSCs",3,1,2.228940333333334,6.686821,0.4998401096732527,0.5213055330575571,1,0.5,0.022297,0.044594,0.1258272041534207,0.08755716033166246,0.04969868261233334,0.298192095674,0.06289348352392458,0.04564403213970329),
    @("MuSCs","Ntf3","Ntrk1","ECs",2,1,0.550853,1.101706,0.1235288445438454,0.08588915922868416,1,0.3333333333333333,0.0009183333333333334,0.002755,0.005182370534790543,0.005409247358696911,0.0005058666716666667,0.00303520003,0.0006401722441607459,0.0004645957076984583),
    @("MuSCs","Ntf3","Ntrk1","FAPs",2,1,0.550853,1.101706,0.1235288445438454,0.08588915922868416,2,0.6666666666666666,0.153988,0.461964,0.8689904253117887,0.9070335923096405,0.08482475176399999,0.508948510584,0.10734538315843,0.07790435263564811),
    @("MuSCs","Ntf3","Ntrk1","MuSCs",2,1,0.550853,1.101706,0.1235288445438454,0.08588915922868416,1,0.5,0.022297,0.044594,0.1258272041534207,0.08755716033166246,0.012282369341,0.04912947736400001,0.01554328914125461,0.007520210885337586)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $row.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}
